$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 219 values that changed
$ws.Range("B219").Value = 1.55
$ws.Range("C219").Value = 1.04
$ws.Range("E219").Value = -0.17
$ws.Range("I219").Value = 2.55
$ws.Range("J219").Value = 8.32
$ws.Range("L219").Value = 6.76
$ws.Range("P219").Value = 5.9

# Row 220
$ws.Range("A220").NumberFormat = "@"
$ws.Range("A220").Value = "03-11-2021"
$ws.Range("A220").ClearFormats()
$ws.Range("B220").Value = 1.61
$ws.Range("C220").Value = 1.07
$ws.Range("D220").Value = 0.08
$ws.Range("E220").Value = -0.17
$ws.Range("F220").Value = 2.5
$ws.Range("G220").Value = 3.62
$ws.Range("H220").Value = 2.98
$ws.Range("I220").Value = 2.51
$ws.Range("J220").Value = 8.300000000000001
$ws.Range("K220").Value = 1.93
$ws.Range("L220").Value = 6.7
$ws.Range("M220").Value = 11.84
$ws.Range("O220").Value = 7.54
$ws.Range("P220").Value = 5.88

# Row 221
$ws.Range("A221").NumberFormat = "@"
$ws.Range("A221").Value = "04-11-2021"
$ws.Range("A221").ClearFormats()
$ws.Range("B221").Value = 1.53
$ws.Range("C221").Value = 0.9399999999999999
$ws.Range("D221").Value = 0.07000000000000001
$ws.Range("E221").Value = -0.22
$ws.Range("F221").Value = 2.45
$ws.Range("H221").Value = 2.9
$ws.Range("I221").Value = 2.65
$ws.Range("K221").Value = 1.95
$ws.Range("L221").Value = 6.55
$ws.Range("M221").Value = 11.97
$ws.Range("O221").Value = 7.45
$ws.Range("P221").Value = 5.85

# Row 222
$ws.Range("A222").NumberFormat = "@"
$ws.Range("A222").Value = "05-11-2021"
$ws.Range("A222").ClearFormats()
$ws.Range("B222").Value = 1.45
$ws.Range("C222").Value = 0.84
$ws.Range("D222").Value = 0.06
$ws.Range("E222").Value = -0.28
$ws.Range("F222").Value = 2.38
$ws.Range("G222").Value = 3.56
$ws.Range("H222").Value = 2.85
$ws.Range("I222").Value = 2.64
$ws.Range("K222").Value = 1.93
$ws.Range("L222").Value = 6.41
$ws.Range("M222").Value = 11.81
$ws.Range("O222").Value = 7.3
$ws.Range("P222").Value = 5.79

# Row 223
$ws.Range("A223").NumberFormat = "@"
$ws.Range("A223").Value = "08-11-2021"
$ws.Range("A223").ClearFormats()
$ws.Range("B223").Value = 1.48
$ws.Range("C223").Value = 0.87
$ws.Range("D223").Value = 0.06
$ws.Range("E223").Value = -0.26
$ws.Range("F223").Value = 2.37
$ws.Range("G223").Value = 3.51
$ws.Range("H223").Value = 2.91
$ws.Range("I223").Value = 2.61
$ws.Range("J223").Value = 8.210000000000001
$ws.Range("K223").Value = 1.9
$ws.Range("L223").Value = 6.32
$ws.Range("M223").Value = 11.95
$ws.Range("O223").Value = 7.25
$ws.Range("P223").Value = 5.8
